$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Click on menu items", "Fail"),
    @("Click on menu items", "Fail"),
    @("Click on menu items", "Fail"),
    @("Send and Delete Message", "Fail"),
    @("Send and Delete Message", "Fail"),
    @("Send and Delete Message", "Fail"),
    @("Send and Delete Message", "Fail"),
    @("Finance Functionality", "Fail"),
    @("Click on menu items", "Fail"),
    @("Click on menu items", "Fail"),
    @("Click on menu items", "Fail"),
    @("Click on menu items", "Fail")
)

$row = 19
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}
